# Generate Report for Handback
# ------------------------------------------------------------------
# The localization round-trip for both target files (zh-cn, de-de) has
# come back "in sync" with en-US, so:
#   * the Status column flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears
#     (Overview sheet + each language sheet),
#   * each language sheet grows a "Latest Target File" / "Latest
#     Handback File" pair of hyperlinked filenames (columns F/G),
#   * the "Latest Handback DateTime" column (H) records when the
#     handback happened.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusNew
$ov.Range("C2").Value = $statusNew
$ov.Range("B3").Value = $statusNew
$ov.Range("C3").Value = $statusNew

# ---- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zh.Range("F2").Value = "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.md"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/292cea940131324afa5cd890b2bf01a3c76218b9/e2e/9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.md", "", "", "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.md") | Out-Null

$zh.Range("G2").Value = "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.6e320472c60673e0b5fda36a56ef48de7758ab26.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec196e36830a0c47eb2287da41eab2c38a1dc1bd/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.6e320472c60673e0b5fda36a56ef48de7758ab26.zh-cn.xlf", "", "", "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.6e320472c60673e0b5fda36a56ef48de7758ab26.zh-cn.xlf") | Out-Null

$zh.Range("F3").Value = "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.md"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/292cea940131324afa5cd890b2bf01a3c76218b9/e2e/c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.md", "", "", "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.md") | Out-Null

$zh.Range("G3").Value = "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.ac4532f7f7e79746e5719f66a3f304ca0ad767ad.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec196e36830a0c47eb2287da41eab2c38a1dc1bd/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.ac4532f7f7e79746e5719f66a3f304ca0ad767ad.zh-cn.xlf", "", "", "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.ac4532f7f7e79746e5719f66a3f304ca0ad767ad.zh-cn.xlf") | Out-Null

$zh.Range("H2").Value = "2016-03-22 07:26:03"
$zh.Range("H3").Value = "2016-03-22 07:26:03"

# ---- de-de sheet ------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$de.Range("F2").Value = "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.md"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/292cea940131324afa5cd890b2bf01a3c76218b9/e2e/9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.md", "", "", "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.md") | Out-Null

$de.Range("G2").Value = "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.6e320472c60673e0b5fda36a56ef48de7758ab26.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/066b834422ce99d87bc3358a3aa2e8f4fba1a271/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.6e320472c60673e0b5fda36a56ef48de7758ab26.de-de.xlf", "", "", "9cf4c740-4aab-4e4d-9fce-d2906f0ab6d4.6e320472c60673e0b5fda36a56ef48de7758ab26.de-de.xlf") | Out-Null

$de.Range("F3").Value = "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.md"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/292cea940131324afa5cd890b2bf01a3c76218b9/e2e/c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.md", "", "", "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.md") | Out-Null

$de.Range("G3").Value = "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.ac4532f7f7e79746e5719f66a3f304ca0ad767ad.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/066b834422ce99d87bc3358a3aa2e8f4fba1a271/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.ac4532f7f7e79746e5719f66a3f304ca0ad767ad.de-de.xlf", "", "", "c1fa5b7e-d1ac-4d67-a387-6ffce1b1414b.ac4532f7f7e79746e5719f66a3f304ca0ad767ad.de-de.xlf") | Out-Null

$de.Range("H2").Value = "2016-03-22 07:26:17"
$de.Range("H3").Value = "2016-03-22 07:26:17"
